{"js": "const body = context.document.body;\n\n// Ordered list of [searchText, [replacement, replacement, ...]] -\n// when searchText occurs more than once in the document, replacements\n// are applied to each match in document order.\nconst replacements = [\n  [\"N = 80,096\", [\"N = 30,424\"]],\n  [\"96.2 (89.6, 102.4)\", [\"96.2 (89.6, 102.5)\"]],\n  [\"79,812 (100)\", [\"30,320 (100)\"]],\n  [\"79,525 (99)\", [\"30,206 (99)\"]],\n  [\"[1.83,319]\", [\"[2.17,323]\"]],\n  [\"20,101 (25)\", [\"7,637 (25)\"]],\n  [\"(319,464]\", [\"(323,467]\"]],\n  [\"20,104 (25)\", [\"7,612 (25)\"]],\n  [\"(464,641]\", [\"(467,641]\"]],\n  [\"19,991 (25)\", [\"7,615 (25)\", \"7,585 (25)\"]],\n  [\"(641,2.39e+03]\", [\"(641,2.14e+03]\"]],\n  [\"19,900 (25)\", [\"7,560 (25)\"]],\n  [\"[0,115]\", [\"[0,117]\"]],\n  [\"19,994 (25)\", [\"7,603 (25)\"]],\n  [\"(115,234]\", [\"(117,236]\"]],\n  [\"20,142 (25)\", [\"7,639 (25)\"]],\n  [\"(234,404]\", [\"(236,410]\"]],\n  [\"20,059 (25)\", [\"7,604 (25)\"]],\n  [\"(404,2.47e+03]\", [\"(410,2.47e+03]\"]],\n  [\"19,901 (25)\", [\"7,578 (25)\"]],\n  [\"[0.978,372]\", [\"[1.02,375]\"]],\n  [\"20,210 (25)\", [\"7,669 (25)\"]],\n  [\"(372,586]\", [\"(375,589]\"]],\n  [\"20,058 (25)\", [\"7,619 (25)\"]],\n  [\"(586,853]\", [\"(589,857]\"]],\n  [\"(853,2.49e+03]\", [\"(857,2.49e+03]\"]],\n  [\"19,837 (25)\", [\"7,551 (25)\"]],\n  [\"462.8 (318.8, 639.8)\", [\"466.3 (322.8, 640.1)\"]],\n  [\"707.6 (515.6, 936.2)\", [\"304.8 (199.3, 441.6)\"]],\n  [\"304.7 (198.3, 443.2)\", [\"138.2 (82.3, 220.6)\"]],\n  [\"233.0 (115.0, 402.5)\", [\"235.6 (117.2, 408.7)\"]],\n  [\"583.4 (370.0, 849.4)\", [\"587.1 (373.6, 854.2)\"]],\n  [\"1,802 (2.2)\", [\"694 (2.3)\"]],\n  [\"1,163 (1.5)\", [\"447 (1.5)\"]],\n  [\"63.3 (56.1, 68.4)\", [\"63.7 (56.6, 68.6)\"]],\n  [\"45,954 (57)\", [\"17,484 (57)\"]],\n  [\"34,142 (43)\", [\"12,940 (43)\"]],\n  [\"6,319 (7.9)\", [\"2,395 (7.9)\"]],\n  [\"20,152 (25)\", [\"7,551 (25)\"]],\n  [\"18,911 (24)\", [\"7,153 (24)\"]],\n  [\"34,714 (43)\", [\"13,325 (44)\"]],\n  [\"0.1 (0.0, 0.6)\", [\"0.1 (0.0, 0.7)\"]],\n  [\"29.8 (16.8, 50.3)\", [\"29.9 (16.8, 50.6)\"]],\n  [\"19.9 (6.0, 41.7)\", [\"19.9 (6.0, 42.0)\"]],\n];\n\nfor (const [searchText, newTexts] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== newTexts.length) {\n    throw new Error(\n      `Expected ${newTexts.length} match(es) for ${JSON.stringify(searchText)}, found ${results.items.length}`\n    );\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newTexts[i], \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Ordered list of (searchText, replacement list) pairs. When searchText\n# occurs more than once in the document, each occurrence (in document\n# order) is replaced with the corresponding entry from its list.\n$replacements = @(\n    , @(\"N = 80,096\", @(\"N = 30,424\"))\n    , @(\"96.2 (89.6, 102.4)\", @(\"96.2 (89.6, 102.5)\"))\n    , @(\"79,812 (100)\", @(\"30,320 (100)\"))\n    , @(\"79,525 (99)\", @(\"30,206 (99)\"))\n    , @(\"[1.83,319]\", @(\"[2.17,323]\"))\n    , @(\"20,101 (25)\", @(\"7,637 (25)\"))\n    , @(\"(319,464]\", @(\"(323,467]\"))\n    , @(\"20,104 (25)\", @(\"7,612 (25)\"))\n    , @(\"(464,641]\", @(\"(467,641]\"))\n    , @(\"19,991 (25)\", @(\"7,615 (25)\", \"7,585 (25)\"))\n    , @(\"(641,2.39e+03]\", @(\"(641,2.14e+03]\"))\n    , @(\"19,900 (25)\", @(\"7,560 (25)\"))\n    , @(\"[0,115]\", @(\"[0,117]\"))\n    , @(\"19,994 (25)\", @(\"7,603 (25)\"))\n    , @(\"(115,234]\", @(\"(117,236]\"))\n    , @(\"20,142 (25)\", @(\"7,639 (25)\"))\n    , @(\"(234,404]\", @(\"(236,410]\"))\n    , @(\"20,059 (25)\", @(\"7,604 (25)\"))\n    , @(\"(404,2.47e+03]\", @(\"(410,2.47e+03]\"))\n    , @(\"19,901 (25)\", @(\"7,578 (25)\"))\n    , @(\"[0.978,372]\", @(\"[1.02,375]\"))\n    , @(\"20,210 (25)\", @(\"7,669 (25)\"))\n    , @(\"(372,586]\", @(\"(375,589]\"))\n    , @(\"20,058 (25)\", @(\"7,619 (25)\"))\n    , @(\"(586,853]\", @(\"(589,857]\"))\n    , @(\"(853,2.49e+03]\", @(\"(857,2.49e+03]\"))\n    , @(\"19,837 (25)\", @(\"7,551 (25)\"))\n    , @(\"462.8 (318.8, 639.8)\", @(\"466.3 (322.8, 640.1)\"))\n    , @(\"707.6 (515.6, 936.2)\", @(\"304.8 (199.3, 441.6)\"))\n    , @(\"304.7 (198.3, 443.2)\", @(\"138.2 (82.3, 220.6)\"))\n    , @(\"233.0 (115.0, 402.5)\", @(\"235.6 (117.2, 408.7)\"))\n    , @(\"583.4 (370.0, 849.4)\", @(\"587.1 (373.6, 854.2)\"))\n    , @(\"1,802 (2.2)\", @(\"694 (2.3)\"))\n    , @(\"1,163 (1.5)\", @(\"447 (1.5)\"))\n    , @(\"63.3 (56.1, 68.4)\", @(\"63.7 (56.6, 68.6)\"))\n    , @(\"45,954 (57)\", @(\"17,484 (57)\"))\n    , @(\"34,142 (43)\", @(\"12,940 (43)\"))\n    , @(\"6,319 (7.9)\", @(\"2,395 (7.9)\"))\n    , @(\"20,152 (25)\", @(\"7,551 (25)\"))\n    , @(\"18,911 (24)\", @(\"7,153 (24)\"))\n    , @(\"34,714 (43)\", @(\"13,325 (44)\"))\n    , @(\"0.1 (0.0, 0.6)\", @(\"0.1 (0.0, 0.7)\"))\n    , @(\"29.8 (16.8, 50.3)\", @(\"29.9 (16.8, 50.6)\"))\n    , @(\"19.9 (6.0, 41.7)\", @(\"19.9 (6.0, 42.0)\"))\n)\n\nforeach ($pair in $replacements) {\n    $searchText = $pair[0]\n    $newTexts = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.MatchWholeWord = $false\n    $i = 0\n    while ($rng.Find.Execute($searchText)) {\n        if ($i -ge $newTexts.Length) {\n            throw \"More matches than expected replacements for $searchText\"\n        }\n        $rng.Text = $newTexts[$i]\n        $rng.Collapse(0)\n        $i++\n    }\n    if ($i -ne $newTexts.Length) {\n        throw \"Expected $($newTexts.Length) match(es) for $searchText, found $i\"\n    }\n}\n"}
